# Apply "Add data for 2022-06-02" update to the carjacking by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-05-25"

# Update the header label in I1 (shared string) to match the new date.
$ws.Range("I1").Value = "2022 (through 05-25)"

# Update individual month figures for the "2022 (through ...)" column (I)
# and year 2021 column (H) that changed with the new day's data.
$ws.Range("I5").Value = 116   # May
$ws.Range("I6").Value = 93    # June

$ws.Range("H12").Value = 202  # November

# Update the Total row (14)
$ws.Range("H14").Value = 1850
$ws.Range("I14").Value = 644
